$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price cells are plain text in the source data (e.g. "35.315.48" or
# "12.60"); Excel auto-converts anything that parses as a number when you assign
# .Value directly, which both changes the cell type and can drop trailing zeros.
# Force text via NumberFormat "@" around the assignment, then restore the "Normal"
# style so the cell-level style index matches the original (unstyled) cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.315.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.911.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.717"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "253.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.78%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.67"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.35%  "
$ws.Range("E9").Value = "  +1.63%  "
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0748"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0991"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.188.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.717"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.901.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.327.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0847"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "243.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.54%  "
$ws.Range("E23").Value = "  +4.33%  "
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.22%  "
$ws.Range("E26").Value = "  +2.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.132"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.126.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +19.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +14.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +22.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0582"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.20"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.914"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0218"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.22%  "
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0646"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.339.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("E48").Value = "  +2.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "12.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +15.44%  "
